$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New period order (ascending) for YERALDIN's 12 rows, then MARIELE's single
# row (period 2103), then YERALDIN's final row (period 2104) appended last.
$periods = @("2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 35112
    $ws.Cells.Item($r, 7).Value = 877803
}

# Row 28: Mariele del Mar Piña Pajaro, periodo 2103
$ws.Cells.Item(28, 2).Value = "CC"
$ws.Cells.Item(28, 3).Value = "1143401993"
$ws.Cells.Item(28, 4).Value = "MARIELE DEL MAR PIÑA PAJARO"
$ws.Cells.Item(28, 5).Value = "2103"
$ws.Cells.Item(28, 6).Value = 9691
$ws.Cells.Item(28, 7).Value = 908526

# Row 29: Yeraldin's final entry, periodo 2104
$ws.Cells.Item(29, 2).Value = "CC"
$ws.Cells.Item(29, 3).Value = "1047438450"
$ws.Cells.Item(29, 4).Value = "YERALDIN PAOLA MERCADO OLIVERO"
$ws.Cells.Item(29, 5).Value = "2104"
$ws.Cells.Item(29, 6).Value = 30430
$ws.Cells.Item(29, 7).Value = 877803
